# Auto-generated Excel COM-interop script to apply the Cactuar_Profits edits
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 17950
$ws.Range("J7").Value = 11000
$ws.Range("L7").Value = 11000
$ws.Range("N7").Value = -11224
# Row 14
$ws.Range("H14").Value = 17950
$ws.Range("J14").Value = 11000
$ws.Range("L14").Value = 11000
$ws.Range("N14").Value = -11382
# Row 33
$ws.Range("H33").Value = 490.13333
$ws.Range("I33").Value = 537.5
$ws.Range("K33").Value = 537.5
$ws.Range("M33").Value = -308.5
# Row 48
$ws.Range("H48").Value = 6124.9165
$ws.Range("J48").Value = 7799.8
$ws.Range("L48").Value = 23399.4
$ws.Range("N48").Value = -23983.4
# Row 56
$ws.Range("H56").Value = 6124.9165
$ws.Range("J56").Value = 7799.8
$ws.Range("L56").Value = 23399.4
$ws.Range("N56").Value = -24467.4
# Row 86
$ws.Range("H86").Value = 100001070
$ws.Range("I86").Value = 142858260
$ws.Range("K86").Value = 142858260
$ws.Range("M86").Value = -142857137
# Row 89
$ws.Range("H89").Value = 100001070
$ws.Range("I89").Value = 142858260
$ws.Range("K89").Value = 714291300
$ws.Range("M89").Value = -714285684
# Row 98
$ws.Range("H98").Value = 1000
$ws.Range("I98").Value = 1000
$ws.Range("K98").Value = 1000
$ws.Range("M98").Value = 498
# Row 106
$ws.Range("H106").Value = 2068.6
$ws.Range("I106").Value = 1448
$ws.Range("K106").Value = 1448
$ws.Range("M106").Value = -817
# Row 116
$ws.Range("H116").Value = 24641718
$ws.Range("I116").Value = 25761438
$ws.Range("K116").Value = 25761438
$ws.Range("M116").Value = -25757996
# Row 122
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
# Row 125
$ws.Range("H125").Value = 2321.875
$ws.Range("I125").Value = 1595.8334
$ws.Range("K125").Value = 14362.5006
$ws.Range("M125").Value = -11902.5006
# Row 132
$ws.Range("H132").Value = 6556.3105
$ws.Range("I132").Value = 1888.381
$ws.Range("K132").Value = 5665.143
$ws.Range("M132").Value = -3135.143
# Row 135
$ws.Range("H135").Value = 6867.6
$ws.Range("I135").Value = 2008.4
$ws.Range("J135").Value = 11726.8
$ws.Range("K135").Value = 18075.6
$ws.Range("L135").Value = 105541.2
$ws.Range("M135").Value = -15540.6
$ws.Range("N135").Value = -110611.2

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 922.82355
$ws.Range("I2").Value = 924.25
$ws.Range("K2").Value = 924.25
$ws.Range("M2").Value = -811.25
# Row 116
$ws.Range("H116").Value = 922.82355
$ws.Range("I116").Value = 924.25
$ws.Range("K116").Value = 924.25
$ws.Range("M116").Value = 1369.75
# Row 122
$ws.Range("H122").Value = 7428.75
$ws.Range("I122").Value = 6905.5557
$ws.Range("K122").Value = 20716.6671
$ws.Range("M122").Value = -18266.6671

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 922.82355
$ws.Range("I3").Value = 924.25
$ws.Range("K3").Value = 924.25
$ws.Range("M3").Value = -810.25
# Row 134
$ws.Range("H134").Value = 2564.0952
$ws.Range("I134").Value = 2182
$ws.Range("K134").Value = 6546
$ws.Range("M134").Value = -4011

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2192.5
$ws.Range("I22").Value = 2773
$ws.Range("K22").Value = 2773
$ws.Range("M22").Value = -2423
# Row 31
$ws.Range("H31").Value = 13701742
$ws.Range("I31").Value = 15627720
$ws.Range("J31").Value = 5896.8887
$ws.Range("K31").Value = 15627720
$ws.Range("L31").Value = 5896.8887
$ws.Range("M31").Value = -15627425
$ws.Range("N31").Value = -6486.8887
# Row 34
$ws.Range("H34").Value = 13701742
$ws.Range("I34").Value = 15627720
$ws.Range("J34").Value = 5896.8887
$ws.Range("K34").Value = 15627720
$ws.Range("L34").Value = 5896.8887
$ws.Range("M34").Value = -15627518
$ws.Range("N34").Value = -6300.8887
# Row 103
$ws.Range("H103").Value = 19076.5
$ws.Range("I103").Value = 8810
$ws.Range("K103").Value = 8810
$ws.Range("M103").Value = -7638

$ws = $wb.Worksheets.Item("CUL")
# Row 28
$ws.Range("H28").Value = 3715
$ws.Range("I28").Value = 5030
$ws.Range("K28").Value = 15090
$ws.Range("M28").Value = -14858
# Row 34
$ws.Range("H34").Value = 608032.5
$ws.Range("J34").Value = 2549.6
$ws.Range("L34").Value = 7648.799999999999
$ws.Range("N34").Value = -7816.799999999999
# Row 114
$ws.Range("H114").Value = 1652.3334
$ws.Range("I114").Value = 481
$ws.Range("J114").Value = 1987
$ws.Range("K114").Value = 1443
$ws.Range("L114").Value = 5961
$ws.Range("M114").Value = 1811
$ws.Range("N114").Value = -12469
# Row 121
$ws.Range("H121").Value = 764.2632
$ws.Range("J121").Value = 855.5625
$ws.Range("L121").Value = 2566.6875
$ws.Range("N121").Value = -5186.6875
# Row 131
$ws.Range("H131").Value = 3065.0833
$ws.Range("J131").Value = 3191.375
$ws.Range("L131").Value = 9574.125
$ws.Range("N131").Value = -19654.125

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("I102").Value = 21741604
$ws.Range("J102").Value = 5306.5
$ws.Range("K102").Value = 21741604
$ws.Range("L102").Value = 5306.5
$ws.Range("M102").Value = -21739982
$ws.Range("N102").Value = -8550.5
# Row 126
$ws.Range("H126").Value = 4698
$ws.Range("I126").Value = 5270.154
$ws.Range("J126").Value = 4078.1667
$ws.Range("K126").Value = 15810.462
$ws.Range("L126").Value = 12234.5001
$ws.Range("M126").Value = -13340.462
$ws.Range("N126").Value = -17174.5001
# Row 132
$ws.Range("H132").Value = 337393.66
$ws.Range("J132").Value = 9425
$ws.Range("L132").Value = 28275
$ws.Range("N132").Value = -33335

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 945.26666
$ws.Range("I22").Value = 827.8
$ws.Range("J22").Value = 1180.2
$ws.Range("K22").Value = 827.8
$ws.Range("L22").Value = 1180.2
$ws.Range("M22").Value = -532.8
$ws.Range("N22").Value = -1770.2
# Row 27
$ws.Range("H27").Value = 945.26666
$ws.Range("I27").Value = 827.8
$ws.Range("J27").Value = 1180.2
$ws.Range("K27").Value = 827.8
$ws.Range("L27").Value = 1180.2
$ws.Range("M27").Value = -720.8
$ws.Range("N27").Value = -1394.2
# Row 118
$ws.Range("H118").Value = 120000
$ws.Range("J118").Value = 120000
$ws.Range("L118").Value = 120000
$ws.Range("N118").Value = -123314
# Row 122
$ws.Range("H122").Value = 12997
$ws.Range("I122").Value = 5002
$ws.Range("J122").Value = 14596
$ws.Range("K122").Value = 15006
$ws.Range("L122").Value = 43788
$ws.Range("M122").Value = -12556
$ws.Range("N122").Value = -48688
# Row 132
$ws.Range("H132").Value = 2731.7
$ws.Range("I132").Value = 2703.9868
$ws.Range("K132").Value = 8111.9604
$ws.Range("M132").Value = -5581.9604
# Row 136
$ws.Range("H136").Value = 4424.956
$ws.Range("I136").Value = 2868.054
$ws.Range("K136").Value = 8604.162
$ws.Range("M136").Value = -6054.162

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 7745.9287
$ws.Range("I100").Value = 649.46155
$ws.Range("K100").Value = 1298.9231
$ws.Range("M100").Value = -757.9231
# Row 110
$ws.Range("H110").Value = 39500
$ws.Range("J110").Value = 39500
$ws.Range("L110").Value = 39500
$ws.Range("N110").Value = -47680
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
# Row 132
$ws.Range("H132").Value = 7251554.5
$ws.Range("I132").Value = 19608982
$ws.Range("K132").Value = 58826946
$ws.Range("M132").Value = -58824416
# Row 136
$ws.Range("H136").Value = 4313.675
$ws.Range("I136").Value = 2987.423
$ws.Range("K136").Value = 8962.269
$ws.Range("M136").Value = -6412.269
